$d = $word.ActiveDocument

# Insert a new paragraph after the last paragraph in the document
# ("Cronograma (CRO): <trello link> ") for the new "Estimativa de
# Tamanho (EST)" link entry.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Seed the new paragraph with its trailing space run first so it is not
# empty when the hyperlink is added (this host anchors Hyperlinks.Add at
# the start of the containing paragraph when it is empty, leaving a stray
# empty run behind; seeding avoids that).
$r = $newPara.Range
$r.Collapse(0)
$r.InsertAfter(" ")

# Add the hyperlink; it lands before the trailing space run.
$hPos = $newPara.Range.Start
$hRange = $d.Range($hPos, $hPos)
$url = "https://docs.google.com/spreadsheets/d/19qxqm7kSXSejTXjsnCBvbB7dTF4EpYtGKVSHF_9FBwI/edit?usp=sharing"
$d.Hyperlinks.Add($hRange, $url, [Type]::Missing, [Type]::Missing, $url, [Type]::Missing)

# Prepend the ": " separator and then the label, each as their own run,
# matching the pattern used by the other entries in the document.
$sepPos = $newPara.Range.Start
$sepRange = $d.Range($sepPos, $sepPos)
$sepRange.InsertBefore(": ")

$labelPos = $newPara.Range.Start
$labelRange = $d.Range($labelPos, $labelPos)
$labelRange.InsertBefore("Estimativa de Tamanho (EST)")
